$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled data-cell format) captured from an untouched cell
# so that forcing text on numeric-looking values does not introduce a new style index.
$styleRef = $ws.Range("D27").Style

# Row 2
$ws.Range("D2").Value = "42.609.65"
$ws.Range("E2").Value = "  -5.75%  "

# Row 3
$ws.Range("D3").Value = "2.213.63"
$ws.Range("E3").Value = "  -6.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.87"
$ws.Range("D5").Style = $styleRef
$ws.Range("E5").Value = "  +2.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.59"
$ws.Range("D6").Style = $styleRef
$ws.Range("E6").Value = "  -9.88%  "

# Row 7
$ws.Range("E7").Value = "  -9.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("E9").Value = "  -9.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.46"
$ws.Range("D10").Style = $styleRef
$ws.Range("E10").Value = "  -11.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.63"
$ws.Range("D11").Style = $styleRef
$ws.Range("E11").Value = "  -3.35%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0821"
$ws.Range("D12").Style = $styleRef
$ws.Range("E12").Value = "  -10.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.63"
$ws.Range("D13").Style = $styleRef
$ws.Range("E13").Value = "  -9.85%  "

# Row 14
$ws.Range("E14").Value = "  -3.99%  "

# Row 15
$ws.Range("D15").Value = "2.547.82"
$ws.Range("E15").Value = "  -6.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("D16").Style = $styleRef
$ws.Range("E16").Value = "  -13.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("D17").Style = $styleRef
$ws.Range("E17").Value = "  -9.36%  "

# Row 18
$ws.Range("D18").Value = "2.195.80"
$ws.Range("E18").Value = "  -6.51%  "

# Row 19
$ws.Range("D19").Value = "42.513.58"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("D20").Style = $styleRef
$ws.Range("E20").Value = "  +5.99%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0954"
$ws.Range("E21").Value = "  -10.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("D22").Style = $styleRef
$ws.Range("E22").Value = "  -11.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.58"
$ws.Range("D23").Style = $styleRef
$ws.Range("E23").Value = "  -10.65%  "

# Row 24
$ws.Range("E24").Value = "  -7.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "233.75"
$ws.Range("D25").Style = $styleRef
$ws.Range("E25").Value = "  -10.05%  "

# Row 26
$ws.Range("E26").Value = "  -6.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = $styleRef
$ws.Range("E28").Value = "  -8.80%  "

# Row 29
$ws.Range("E29").Value = "  -8.27%  "

# Row 30
$ws.Range("E30").Value = "  -14.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.44"
$ws.Range("D31").Style = $styleRef
$ws.Range("E31").Value = "  -8.79%  "

# Row 32
$ws.Range("E32").Value = "  -8.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.84"
$ws.Range("D33").Style = $styleRef
$ws.Range("E33").Value = "  -6.74%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.68"
$ws.Range("D34").Style = $styleRef
$ws.Range("E34").Value = "  -10.52%  "

# Row 35
$ws.Range("E35").Value = "  -7.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.20"
$ws.Range("D36").Style = $styleRef
$ws.Range("E36").Value = "  +7.25%  "

# Row 37
$ws.Range("E37").Value = "  -7.11%  "

# Row 38
$ws.Range("E38").Value = "  +7.56%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.41"
$ws.Range("D39").Style = $styleRef
$ws.Range("E39").Value = "  -8.23%  "

# Row 40
$ws.Range("E40").Value = "  -11.05%  "

# Row 41
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("D41").Style = $styleRef
$ws.Range("E41").Value = "  -10.86%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0318"
$ws.Range("D42").Style = $styleRef
$ws.Range("E42").Value = "  -10.50%  "

# Row 43
$ws.Range("E43").Value = "  +0.33%  "

# Row 44
$ws.Range("D44").Value = "1.772.25"
$ws.Range("E44").Value = "  +8.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.65"
$ws.Range("D45").Style = $styleRef
$ws.Range("E45").Value = "  -12.82%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.88"
$ws.Range("D46").Style = $styleRef
$ws.Range("E46").Value = "  -9.45%  "

# Row 47
$ws.Range("E47").Value = "  -12.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.67"
$ws.Range("D48").Style = $styleRef
$ws.Range("E48").Value = "  -6.33%  "

# Row 49
$ws.Range("E49").Value = "  -4.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.83"
$ws.Range("D50").Style = $styleRef
$ws.Range("E50").Value = "  -14.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.60"
$ws.Range("D51").Style = $styleRef
$ws.Range("E51").Value = "  +56.24%  "
